# daily auto push: 2025-10-12 13:28 UTC
# Append the 2025/10/12 20:00 ranking row (row 96) to the end of the log sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A stores dates as plain text (e.g. "2025/10/12"), not native Excel
# date serials, so force text formatting before assigning the value to stop
# Excel's autodetection from turning the string into a date number. Resetting
# the style back to Normal afterwards keeps the new cell unstyled, matching
# every other data row in the sheet.
$ws.Range("A96").NumberFormat = "@"
$ws.Range("A96").Value = "2025/10/12"
$ws.Range("A96").Style = "Normal"

$ws.Range("B96").Value = "日"
$ws.Range("C96").Value = 20
$ws.Range("D96").Value = 201
